# Javítókulcs 3. feladat - update rubric columns for the new "Teletál dolgozat"
# (Monaco2023.csv / Forma-1 kvalifikáció) task, replacing the previous
# "étlap" task rubric, and clear the now-unused AI:AL sample answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adatok")

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")

# --- New rubric headers (row 1), columns F..W -----------------------------
$headerTexts = @(
    "Létrehoz egy programot kvalifikacio.py néven.",
    "A program hibaüzenet nélkül lefut.",
    "Adatszerkezetet hoz létre a Monaco2023.csv adatainak eltárolására.",
    "A Monaco2023.csv legalább egy sorát beolvassa.",
    "Meghatározza a köridők számát.",
    "A mintának megfelelően kiírja a képernyőre a köridők számát.",
    "Az osztályban tagfüggvényt hoz létre.",
    "A tagfüggvény a köridőt három tizedesjegy pontossággal másodpercre váltva adja vissza.",
    "Meghatározza a Q1-ben legjobb eredményt elért versenyzőt.",
    "A minta szerint kiírta a képernyőre a Q1-ben legjobb eredményt elért versenyző adatait.",
    "Kilistázza a képernyőre a Q2-be jutott versenyzőket, valamint az ott elért köridejüket.",
    "Bekér egy köridőt.",
    "Eltárolja a bekért köridőt.",
    "Létrehozza a hatarfeletti.txt nevű állományt.",
    "Létrehozza statisztikához szükséges objektumot.",
    "Létrehozza a statisztikát a feladatnak megfelelően.",
    "Feltételt állít konstruktőrök köridejeinket számának.",
    "Kilistázza a feltételnek megfelelő konstruktőrök nevét és körideik számát."
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value2 = $headerTexts[$i]
}

# Columns X..AL no longer hold individual rubric items - blank them out.
$ws.Range("X1:AL1").ClearContents()

# --- Row 2: task numbering now 1..18 (one point each, was grouped 0..8) ---
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "2").Value2 = ($i + 1)
}
$ws.Range("X2:AL2").ClearContents()

# --- Row 3: max point value per task is 1 for all 18 tasks ---------------
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value2 = 1
}
$ws.Range("X3:AL3").ClearContents()

# --- Sample students: clear the now unused AI:AL answer columns ----------
$ws.Range("AI4:AL4").ClearContents()
$ws.Range("AI5:AL5").ClearContents()
$ws.Range("AI6:AL6").ClearContents()
$ws.Range("AI7:AL7").ClearContents()

# Recalculate all dependent formulas (AO/AP/AQ per student, stats rows, charts, Jegyek sheet)
$excel.CalculateFullRebuild()

$ws.Activate()
$ws.Range("F4").Select()

Write-Host "Rubric update applied."
